# "delete container, fixing screenshot pdf"
# The source data cell A2 ("NPP" column) held a stray/incorrect code that
# is corrected here, and the previously-left-over cell selection is reset
# to the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell A2 is formatted as Text (quote-prefixed numeric string). Prefixing
# the new value with an apostrophe keeps it a genuine text entry (same
# cell style / quote-prefix flag) instead of Excel reinterpreting it as a
# plain number and dropping the text formatting.
$ws.Range("A2").Value = "'01732301"

# Reset the active selection to A3 (matches the saved cursor position).
$ws.Range("A3").Select()
